$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits at the very
#    end of the document, right after the "}" that closes the
#    "{report.comment}" placeholder).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Rename the table-cell heading "Методика" -> "Пункт требований".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Методика", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Пункт требований", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark immediately after the run we just
#    edited (end of that paragraph, inside the table cell).
#
#    A bookmark range collapsed exactly on the paragraph-end boundary can't
#    be targeted directly, so a one-character placeholder is inserted right
#    after the new text, the bookmark is anchored around that placeholder,
#    and the placeholder is then deleted again - the bookmark collapses
#    back to a zero-length bookmark at the correct spot.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Пункт требований", $false, $false, $false, $false,
                   $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("X")
$markerPos = $rng.Start

$bmRange = $d.Range($markerPos, $markerPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Range($markerPos, $markerPos + 1).Delete() | Out-Null
